# Section 1 Visualization Update
# Adds "Sheet2" - a transposed view of Sheet1's crime-against-women data
# (categories as rows, years as columns) - and updates view/selection state.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# --- Add the new worksheet right after Sheet1 ---------------------------
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# --- Populate Sheet2 with the transposed table ---------------------------
# Row 1: "Year" label followed by the 18 year columns (1998-2015).
# Rows 2-13: each crime category label followed by its value for every year.
$sheet2Data = @(
    @('Year', 1998, 1999, 2000, 2001, 2002, 2003, 2004, 2005, 2006, 2007, 2008, 2009, 2010, 2011, 2012, 2013, 2014, 2015),
    @('Rape', 368, 331, 330, 286, 267, 236, 339, 324, 354, 316, 374, 433, 408, 439, 473, 732, 841, 503),
    @('Kidnapping_Abduction', 1182, 1074, 868, 857, 807, 859, 905, 916, 945, 1089, 1119, 1162, 1290, 1442, 1527, 2230, 2187, 1569),
    @('Dowry_Deaths', 90, 94, 93, 67, 62, 54, 58, 48, 50, 42, 27, 24, 19, 30, 21, 29, 23, 12),
    @('Domestic_Violence', 3602, 3886, 3739, 3667, 3321, 3684, 3955, 4090, 4977, 5827, 6094, 5506, 5600, 6052, 6658, 7812, 5991, 4133),
    @('Molestation', 1210, 1083, 944, 756, 750, 722, 757, 802, 736, 822, 828, 727, 668, 685, 745, 1243, 1352, 1164),
    @('Eveteasing', 139, 172, 119, 111, 104, 92, 164, 104, 138, 120, 122, 114, 110, 93, 93, 77, 173, 73),
    @('Importing', 57, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('SATI', 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('Immoral_Traffic', 8, 40, 47, 61, 57, 74, 33, 59, 78, 44, 52, 41, 46, 46, 44, 76, 45, 35),
    @('Indecent Representation', 2, 1, 0, 0, 4, 14, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @('Dowry Proh.Act', 0, 13, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 84, 53, 36),
    @('Total', 6658, 6694, 6140, 5805, 5373, 5735, 6211, 6343, 7279, 8260, 8616, 8009, 8148, 8815, 9561, 12283, 10665, 7525)
)

for ($r = 0; $r -lt $sheet2Data.Length; $r++) {
    $rowValues = $sheet2Data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws2.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

# --- View/selection state -------------------------------------------------
# Sheet1: selection becomes the whole used range, Sheet1 no longer the
# active/focused tab.
$sheet1.Range("A1:M19").Select()

# Sheet2: becomes the active tab with G7 selected.
$ws2.Range("G7").Select()
$ws2.Activate()
